$d = $word.ActiveDocument

# The three "gestore utenti di ..." sentences that must become
# "gestore dipendenti di ...". Word's actual editing behaviour (typing
# "dipendenti" over a selected "utenti") splits the sentence's run into
# three runs around the edited word, even though they end up with
# identical run properties. We reproduce that by toggling Bold on/off
# on the replaced sub-range, which forces the engine to keep the run
# boundaries instead of silently re-merging them.
$sentences = @(
    "Il sistema deve permettere al gestore utenti di inserire un nuovo dipendente.",
    "Il sistema deve permettere al gestore utenti di modificare le informazioni di un dipendente.",
    "Il sistema deve permettere al gestore utenti di eliminare un dipendente."
)

foreach ($sentence in $sentences) {
    $rng = $d.Content
    $found = $rng.Find.Execute($sentence, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        continue
    }

    $word_rng = $rng.Duplicate
    $word_rng.Find.Execute("utenti", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $word_rng.Text = "dipendenti"
    $word_rng.Bold = 1
    $word_rng.Bold = 0
}
